$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for ""Barbie"" to be shown on Friday.`n"
$ws.Range("D2").Value = "Barbie_was_selected, "

$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been registered, and no movie will be shown on Friday.`n"

$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday has not been made.`n"

$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no agreement was reached on which movie to show on Friday.`n"

$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding the movie to be shown on Friday.`n"

$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision regarding the movie to show on Friday is recorded as a no decision.`n"

$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be selected for Friday.`n"

$ws.Range("C9").Value = "MSG: None`n`nMSG: It seems there was no definitive choice made regarding the movie to be shown on Friday, so I must proceed with the no_decision function.`n"
